$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data cells to match the latest scrape.
# Cells whose new text would otherwise be auto-parsed as a number are
# forced to Text format first so the literal string is preserved.

$ws.Range("D2").Value = "27.651.94"
$ws.Range("E2").Value = "  -4.12%  "
$ws.Range("D3").Value = "1.847.09"
$ws.Range("E3").Value = "  -3.58%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.46"
$ws.Range("E5").Value = "  -3.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4271"
$ws.Range("E7").Value = "  -6.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3647"
$ws.Range("E8").Value = "  -4.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.80"
$ws.Range("E9").Value = "  -3.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07231"
$ws.Range("E10").Value = "  -6.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8983"
$ws.Range("E11").Value = "  -8.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.62"
$ws.Range("E12").Value = "  -7.07%  "
$ws.Range("D13").Value = "1.832.32"
$ws.Range("E13").Value = "  -5.83%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.374"
$ws.Range("E14").Value = "  -5.62%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.571"
$ws.Range("E15").Value = "  -5.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06853"
$ws.Range("E16").Value = "  -2.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "77.75"
$ws.Range("E18").Value = "  -7.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008907"
$ws.Range("E19").Value = "  -6.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.42"
$ws.Range("E21").Value = "  -7.32%  "
$ws.Range("D22").Value = "27.620.52"
$ws.Range("E22").Value = "  -4.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.966"
$ws.Range("E23").Value = "  -6.91%  "
$ws.Range("E24").Value = "  -3.84%  "
$ws.Range("D25").Value = "2.051.68"
$ws.Range("E25").Value = "  -4.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.048"
$ws.Range("E26").Value = "  -0.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.32"
$ws.Range("E27").Value = "  -2.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.28"
$ws.Range("E28").Value = "  -4.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.338"
$ws.Range("E29").Value = "  -4.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "111.63"
$ws.Range("E30").Value = "  -5.29%  "
$ws.Range("E31").Value = "  -3.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08927"
$ws.Range("E32").Value = "  -4.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7829"
$ws.Range("E33").Value = "  -10.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.518"
$ws.Range("E34").Value = "  -11.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.876"
$ws.Range("E35").Value = "  -5.03%  "
$ws.Range("E36").Value = "  -12.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.001"
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05440"
$ws.Range("E38").Value = "  -4.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.097"
$ws.Range("E39").Value = "  -4.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.985"
$ws.Range("E40").Value = "  -2.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01925"
$ws.Range("E41").Value = "  -5.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5053"
$ws.Range("E42").Value = "  -8.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.795"
$ws.Range("E43").Value = "  -9.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1637"
$ws.Range("E44").Value = "  -6.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.299"
$ws.Range("E45").Value = "  -11.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06630"
$ws.Range("E46").Value = "  -4.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "106.59"
$ws.Range("E47").Value = "  -3.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4727"
$ws.Range("E48").Value = "  -8.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.38"
$ws.Range("E49").Value = "  -7.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9999"
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.645"
$ws.Range("E51").Value = "  -6.75%  "
